$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H57").Value = 10825
$ws.Range("J57").Value = 10825
$ws.Range("L57").Value = 32475
$ws.Range("N57").Value = -33473
$ws.Range("H106").Value = 23288.53
$ws.Range("I106").Value = 26636.072
$ws.Range("J106").Value = 7666.6665
$ws.Range("K106").Value = 26636.072
$ws.Range("L106").Value = 7666.6665
$ws.Range("M106").Value = -26005.072
$ws.Range("N106").Value = -8928.666499999999
$ws.Range("H112").Value = 2106.818
$ws.Range("J112").Value = 2106.818
$ws.Range("L112").Value = 6320.454000000001
$ws.Range("N112").Value = -8536.454000000002
$ws.Range("H137").Value = 5268.32
$ws.Range("I137").Value = 640.3333
$ws.Range("J137").Value = 6729.7896
$ws.Range("K137").Value = 1920.9999
$ws.Range("L137").Value = 20189.3688
$ws.Range("M137").Value = 629.0001
$ws.Range("N137").Value = -25289.3688

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2669.111
$ws.Range("I61").Value = 2636.3
$ws.Range("J61").Value = 2833.1667
$ws.Range("K61").Value = 2636.3
$ws.Range("L61").Value = 2833.1667
$ws.Range("M61").Value = -2424.3
$ws.Range("N61").Value = -3257.1667
$ws.Range("H74").Value = 103653.8
$ws.Range("I74").Value = 171101.83
$ws.Range("J74").Value = 2481.75
$ws.Range("K74").Value = 171101.83
$ws.Range("L74").Value = 2481.75
$ws.Range("M74").Value = -170227.83
$ws.Range("N74").Value = -4229.75
$ws.Range("H77").Value = 103653.8
$ws.Range("I77").Value = 171101.83
$ws.Range("J77").Value = 2481.75
$ws.Range("K77").Value = 855509.1499999999
$ws.Range("L77").Value = 12408.75
$ws.Range("M77").Value = -851141.1499999999
$ws.Range("N77").Value = -21144.75
$ws.Range("H132").Value = 196541.44
$ws.Range("I132").Value = 30877.107
$ws.Range("J132").Value = 503020.44
$ws.Range("K132").Value = 92631.321
$ws.Range("L132").Value = 1509061.32
$ws.Range("M132").Value = -90101.321
$ws.Range("N132").Value = -1514121.32
$ws.Range("H136").Value = 2669.111
$ws.Range("I136").Value = 2636.3
$ws.Range("J136").Value = 2833.1667
$ws.Range("K136").Value = 7908.900000000001
$ws.Range("L136").Value = 8499.500100000001
$ws.Range("M136").Value = -5358.900000000001
$ws.Range("N136").Value = -13599.5001
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("H141").Value = 96685.8
$ws.Range("J141").Value = 96685.8
$ws.Range("L141").Value = 96685.8
$ws.Range("N141").Value = -107045.8
$ws.Range("N140").ClearContents()

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 22047.8
$ws.Range("J81").Value = 22047.8
$ws.Range("L81").Value = 22047.8
$ws.Range("N81").Value = -24169.8
$ws.Range("H84").Value = 22047.8
$ws.Range("J84").Value = 22047.8
$ws.Range("L84").Value = 66143.39999999999
$ws.Range("N84").Value = -76751.39999999999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4650.8184
$ws.Range("I31").Value = 1302.8158
$ws.Range("J31").Value = 12134.588
$ws.Range("K31").Value = 1302.8158
$ws.Range("L31").Value = 12134.588
$ws.Range("M31").Value = -1007.8158
$ws.Range("N31").Value = -12724.588
$ws.Range("H34").Value = 4650.8184
$ws.Range("I34").Value = 1302.8158
$ws.Range("J34").Value = 12134.588
$ws.Range("K34").Value = 1302.8158
$ws.Range("L34").Value = 12134.588
$ws.Range("M34").Value = -1100.8158
$ws.Range("N34").Value = -12538.588
$ws.Range("H58").Value = 1412.3636
$ws.Range("I58").Value = 1160.8
$ws.Range("J58").Value = 1951.4286
$ws.Range("K58").Value = 1160.8
$ws.Range("L58").Value = 1951.4286
$ws.Range("M58").Value = -957.8
$ws.Range("N58").Value = -2357.4286
$ws.Range("H132").Value = 49414.477
$ws.Range("I132").Value = 84911.5
$ws.Range("J132").Value = 2085.111
$ws.Range("K132").Value = 254734.5
$ws.Range("L132").Value = 6255.333
$ws.Range("M132").Value = -252204.5
$ws.Range("N132").Value = -11315.333
$ws.Range("H136").Value = 1412.3636
$ws.Range("I136").Value = 1160.8
$ws.Range("J136").Value = 1951.4286
$ws.Range("K136").Value = 3482.4
$ws.Range("L136").Value = 5854.2858
$ws.Range("M136").Value = -932.3999999999996
$ws.Range("N136").Value = -10954.2858

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 1778
$ws.Range("I70").Value = 556
$ws.Range("J70").Value = 3000
$ws.Range("K70").Value = 1668
$ws.Range("L70").Value = 9000
$ws.Range("M70").Value = -1353
$ws.Range("N70").Value = -9630
$ws.Range("H73").Value = 1778
$ws.Range("I73").Value = 556
$ws.Range("J73").Value = 3000
$ws.Range("K73").Value = 1668
$ws.Range("L73").Value = 9000
$ws.Range("M73").Value = -576
$ws.Range("N73").Value = -11184
$ws.Range("H75").Value = 37041130
$ws.Range("I75").Value = 613
$ws.Range("J75").Value = 55561384
$ws.Range("K75").Value = 1839
$ws.Range("L75").Value = 166684152
$ws.Range("M75").Value = -841
$ws.Range("N75").Value = -166686148
$ws.Range("H78").Value = 37041130
$ws.Range("I78").Value = 613
$ws.Range("J78").Value = 55561384
$ws.Range("K78").Value = 5517
$ws.Range("L78").Value = 500052456
$ws.Range("M78").Value = -525
$ws.Range("N78").Value = -500062440
$ws.Range("H97").Value = 3067.1052
$ws.Range("J97").Value = 3621.7856
$ws.Range("L97").Value = 10865.3568
$ws.Range("N97").Value = -11857.3568
$ws.Range("H103").Value = 1544.1052
$ws.Range("I103").Value = 499.0909
$ws.Range("J103").Value = 2981
$ws.Range("K103").Value = 1497.2727
$ws.Range("L103").Value = 8943
$ws.Range("M103").Value = -618.2727
$ws.Range("N103").Value = -10701
$ws.Range("H122").Value = 651.1786
$ws.Range("I122").Value = 462.95456
$ws.Range("J122").Value = 1341.3334
$ws.Range("K122").Value = 4166.59104
$ws.Range("L122").Value = 12072.0006
$ws.Range("M122").Value = -1716.59104
$ws.Range("N122").Value = -16972.0006
$ws.Range("H129").Value = 19270628
$ws.Range("I129").Value = 514.875
$ws.Range("J129").Value = 27835124
$ws.Range("K129").Value = 1544.625
$ws.Range("L129").Value = 83505372
$ws.Range("M129").Value = 3455.375
$ws.Range("N129").Value = -83515372
$ws.Range("H131").Value = 678.4343
$ws.Range("I131").Value = 410.10526
$ws.Range("J131").Value = 742.1625
$ws.Range("K131").Value = 1230.31578
$ws.Range("L131").Value = 2226.4875
$ws.Range("M131").Value = 3809.68422
$ws.Range("N131").Value = -12306.4875

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 40093.31
$ws.Range("I132").Value = 1409.3043
$ws.Range("J132").Value = 336670.66
$ws.Range("K132").Value = 4227.9129
$ws.Range("L132").Value = 1010011.98
$ws.Range("M132").Value = -1697.9129
$ws.Range("N132").Value = -1015071.98
$ws.Range("H141").Value = 43071.547
$ws.Range("J141").Value = 43071.547
$ws.Range("L141").Value = 43071.547
$ws.Range("N141").Value = -53431.547

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 32055.117
$ws.Range("I40").Value = 1962.4
$ws.Range("J40").Value = 44593.75
$ws.Range("K40").Value = 1962.4
$ws.Range("L40").Value = 44593.75
$ws.Range("M40").Value = -1826.4
$ws.Range("N40").Value = -44865.75
$ws.Range("H68").Value = 10943.909
$ws.Range("I68").Value = 15500
$ws.Range("J68").Value = 2970.75
$ws.Range("K68").Value = 15500
$ws.Range("L68").Value = 2970.75
$ws.Range("M68").Value = -14751
$ws.Range("N68").Value = -4468.75
$ws.Range("H71").Value = 10943.909
$ws.Range("I71").Value = 15500
$ws.Range("J71").Value = 2970.75
$ws.Range("K71").Value = 77500
$ws.Range("L71").Value = 14853.75
$ws.Range("M71").Value = -73756
$ws.Range("N71").Value = -22341.75
$ws.Range("H132").Value = 312061.34
$ws.Range("I132").Value = 81414.234
$ws.Range("J132").Value = 773355.6
$ws.Range("K132").Value = 244242.702
$ws.Range("L132").Value = 2320066.8
$ws.Range("M132").Value = -241712.702
$ws.Range("N132").Value = -2325126.8
$ws.Range("H136").Value = 1252231.8
$ws.Range("I136").Value = 10000004
$ws.Range("J136").Value = 2550
$ws.Range("K136").Value = 30000012
$ws.Range("L136").Value = 7650
$ws.Range("M136").Value = -29997462
$ws.Range("N136").Value = -12750

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 756.6
$ws.Range("I126").Value = 756.6
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 2269.8
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = 200.1999999999998
$ws.Range("H132").Value = 4117.706
$ws.Range("I132").Value = 770.4167
$ws.Range("J132").Value = 12151.2
$ws.Range("K132").Value = 2311.2501
$ws.Range("L132").Value = 36453.60000000001
$ws.Range("M132").Value = 218.7498999999998
$ws.Range("N132").Value = -41513.60000000001
$ws.Range("H140").Value = 46597.1
$ws.Range("J140").Value = 46597.1
$ws.Range("L140").Value = 46597.1
$ws.Range("N140").Value = -56957.1
$ws.Range("H141").Value = 48626.58
$ws.Range("J141").Value = 48626.58
$ws.Range("L141").Value = 48626.58
$ws.Range("N141").Value = -58986.58
$ws.Range("N126").ClearContents()
